# Auto-generated edit script.
# Applies the cell-value changes described by the Gilgamesh_Profits.xlsx
# commit diff across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# All target cells hold literal (non-formula) numeric values, so each
# change is a plain value write; the two rows whose HQ-profit column
# (N) is dropped from the "after" state are cleared instead of zeroed,
# matching the diff (cell removed, not set to 0).

$wb = $excel.ActiveWorkbook

# ======== Sheet: ALC ========
$ws = $wb.Worksheets.Item("ALC")

# -- Row 3 --
$ws.Cells.Item(3, 8).Value = 0  # H3: 30657 -> 0
$ws.Cells.Item(3, 10).Value = 0  # J3: 30657 -> 0
$ws.Cells.Item(3, 12).Value = 0  # L3: 30657 -> 0
$ws.Cells.Item(3, 14).ClearContents()  # N3: -30885 -> (removed)

# -- Row 15 --
$ws.Cells.Item(15, 8).Value = 3779.4893  # H15: 3317.551 -> 3779.4893
$ws.Cells.Item(15, 9).Value = 3779.4893  # I15: 3317.551 -> 3779.4893
$ws.Cells.Item(15, 11).Value = 11338.4679  # K15: 9952.653 -> 11338.4679
$ws.Cells.Item(15, 13).Value = -11169.4679  # M15: -9783.653 -> -11169.4679

# -- Row 64 --
$ws.Cells.Item(64, 8).Value = 20414136  # H64: 21282568 -> 20414136
$ws.Cells.Item(64, 9).Value = 6052.1904  # I64: 6050.3955 -> 6052.1904
$ws.Cells.Item(64, 10).Value = 142862640  # J64: 250005120 -> 142862640
$ws.Cells.Item(64, 11).Value = 6052.1904  # K64: 6050.3955 -> 6052.1904
$ws.Cells.Item(64, 12).Value = 142862640  # L64: 250005120 -> 142862640
$ws.Cells.Item(64, 13).Value = -5804.1904  # M64: -5802.3955 -> -5804.1904
$ws.Cells.Item(64, 14).Value = -142863136  # N64: -250005616 -> -142863136

# -- Row 67 --
$ws.Cells.Item(67, 8).Value = 20414136  # H67: 21282568 -> 20414136
$ws.Cells.Item(67, 9).Value = 6052.1904  # I67: 6050.3955 -> 6052.1904
$ws.Cells.Item(67, 10).Value = 142862640  # J67: 250005120 -> 142862640
$ws.Cells.Item(67, 11).Value = 6052.1904  # K67: 6050.3955 -> 6052.1904
$ws.Cells.Item(67, 12).Value = 142862640  # L67: 250005120 -> 142862640
$ws.Cells.Item(67, 13).Value = -5194.1904  # M67: -5192.3955 -> -5194.1904
$ws.Cells.Item(67, 14).Value = -142864356  # N67: -250006836 -> -142864356

# -- Row 96 --
$ws.Cells.Item(96, 8).Value = 677.5  # H96: 639.9 -> 677.5
$ws.Cells.Item(96, 9).Value = 387.125  # I96: 340.125 -> 387.125
$ws.Cells.Item(96, 11).Value = 1161.375  # K96: 1020.375 -> 1161.375
$ws.Cells.Item(96, 13).Value = 211.625  # M96: 352.625 -> 211.625

# -- Row 98 --
$ws.Cells.Item(98, 8).Value = 3543.7778  # H98: 4049.25 -> 3543.7778
$ws.Cells.Item(98, 9).Value = 3852.5334  # I98: 4052.5334 -> 3852.5334
$ws.Cells.Item(98, 10).Value = 2000  # J98: 4000 -> 2000
$ws.Cells.Item(98, 11).Value = 3852.5334  # K98: 4052.5334 -> 3852.5334
$ws.Cells.Item(98, 12).Value = 2000  # L98: 4000 -> 2000
$ws.Cells.Item(98, 13).Value = -2354.5334  # M98: -2554.5334 -> -2354.5334
$ws.Cells.Item(98, 14).Value = -4996  # N98: -6996 -> -4996

# -- Row 102 --
$ws.Cells.Item(102, 8).Value = 0  # H102: 30657 -> 0
$ws.Cells.Item(102, 10).Value = 0  # J102: 30657 -> 0
$ws.Cells.Item(102, 12).Value = 0  # L102: 30657 -> 0
$ws.Cells.Item(102, 14).ClearContents()  # N102: -37147 -> (removed)

# -- Row 107 --
$ws.Cells.Item(107, 8).Value = 473  # H107: 413.27777 -> 473
$ws.Cells.Item(107, 9).Value = 300.2857  # I107: 296.26666 -> 300.2857
$ws.Cells.Item(107, 10).Value = 1077.5  # J107: 998.3333 -> 1077.5
$ws.Cells.Item(107, 11).Value = 300.2857  # K107: 296.26666 -> 300.2857
$ws.Cells.Item(107, 12).Value = 1077.5  # L107: 998.3333 -> 1077.5
$ws.Cells.Item(107, 13).Value = 1619.7143  # M107: 1623.73334 -> 1619.7143
$ws.Cells.Item(107, 14).Value = -4917.5  # N107: -4838.3333 -> -4917.5

# -- Row 122 --
$ws.Cells.Item(122, 8).Value = 3543.7778  # H122: 4049.25 -> 3543.7778
$ws.Cells.Item(122, 9).Value = 3852.5334  # I122: 4052.5334 -> 3852.5334
$ws.Cells.Item(122, 10).Value = 2000  # J122: 4000 -> 2000
$ws.Cells.Item(122, 11).Value = 11557.6002  # K122: 12157.6002 -> 11557.6002
$ws.Cells.Item(122, 12).Value = 6000  # L122: 12000 -> 6000
$ws.Cells.Item(122, 13).Value = -9107.600199999999  # M122: -9707.600199999999 -> -9107.600199999999
$ws.Cells.Item(122, 14).Value = -10900  # N122: -16900 -> -10900

# -- Row 132 --
$ws.Cells.Item(132, 8).Value = 5901.107  # H132: 6569.56 -> 5901.107
$ws.Cells.Item(132, 9).Value = 5934.4814  # I132: 6634.9585 -> 5934.4814
$ws.Cells.Item(132, 11).Value = 17803.4442  # K132: 19904.8755 -> 17803.4442
$ws.Cells.Item(132, 13).Value = -15273.4442  # M132: -17374.8755 -> -15273.4442

# -- Row 138 --
$ws.Cells.Item(138, 8).Value = 2689.258  # H138: 2713.476 -> 2689.258
$ws.Cells.Item(138, 9).Value = 2411.4333  # I138: 2516.6897 -> 2411.4333
$ws.Cells.Item(138, 10).Value = 2949.7188  # J138: 2881.3235 -> 2949.7188
$ws.Cells.Item(138, 11).Value = 7234.2999  # K138: 7550.0691 -> 7234.2999
$ws.Cells.Item(138, 12).Value = 8849.1564  # L138: 8643.970499999999 -> 8849.1564
$ws.Cells.Item(138, 13).Value = -2094.2999  # M138: -2410.0691 -> -2094.2999
$ws.Cells.Item(138, 14).Value = -19129.1564  # N138: -18923.9705 -> -19129.1564


# ======== Sheet: ARM ========
$ws = $wb.Worksheets.Item("ARM")

# -- Row 61 --
$ws.Cells.Item(61, 8).Value = 2800.3076  # H61: 2919.739 -> 2800.3076
$ws.Cells.Item(61, 9).Value = 2514.6667  # I61: 2619.6667 -> 2514.6667
$ws.Cells.Item(61, 11).Value = 2514.6667  # K61: 2619.6667 -> 2514.6667
$ws.Cells.Item(61, 13).Value = -2302.6667  # M61: -2407.6667 -> -2302.6667

# -- Row 136 --
$ws.Cells.Item(136, 8).Value = 2800.3076  # H136: 2919.739 -> 2800.3076
$ws.Cells.Item(136, 9).Value = 2514.6667  # I136: 2619.6667 -> 2514.6667
$ws.Cells.Item(136, 11).Value = 7544.000100000001  # K136: 7859.000100000001 -> 7544.000100000001
$ws.Cells.Item(136, 13).Value = -4994.000100000001  # M136: -5309.000100000001 -> -4994.000100000001


# ======== Sheet: BSM ========
$ws = $wb.Worksheets.Item("BSM")

# -- Row 29 --
$ws.Cells.Item(29, 8).Value = 7450  # H29: 6800 -> 7450
$ws.Cells.Item(29, 9).Value = 7450  # I29: 6800 -> 7450
$ws.Cells.Item(29, 11).Value = 7450  # K29: 6800 -> 7450
$ws.Cells.Item(29, 13).Value = -7161  # M29: -6511 -> -7161

# -- Row 94 --
$ws.Cells.Item(94, 8).Value = 80001180  # H94: 83334800 -> 80001180
$ws.Cells.Item(94, 9).Value = 133333590  # I94: 142857650 -> 133333590
$ws.Cells.Item(94, 10).Value = 2560.4  # J94: 2810.4 -> 2560.4
$ws.Cells.Item(94, 11).Value = 133333590  # K94: 142857650 -> 133333590
$ws.Cells.Item(94, 12).Value = 2560.4  # L94: 2810.4 -> 2560.4
$ws.Cells.Item(94, 13).Value = -133333139  # M94: -142857199 -> -133333139
$ws.Cells.Item(94, 14).Value = -3462.4  # N94: -3712.4 -> -3462.4

# -- Row 99 --
$ws.Cells.Item(99, 8).Value = 4726.913  # H99: 4640.9585 -> 4726.913
$ws.Cells.Item(99, 9).Value = 3199.7144  # I99: 3164 -> 3199.7144
$ws.Cells.Item(99, 11).Value = 3199.7144  # K99: 3164 -> 3199.7144
$ws.Cells.Item(99, 13).Value = -1701.7144  # M99: -1666 -> -1701.7144


# ======== Sheet: CRP ========
$ws = $wb.Worksheets.Item("CRP")

# -- Row 31 --
$ws.Cells.Item(31, 8).Value = 3571.7317  # H31: 3572.2195 -> 3571.7317
$ws.Cells.Item(31, 9).Value = 2619.4243  # I31: 2620.0303 -> 2619.4243
$ws.Cells.Item(31, 11).Value = 2619.4243  # K31: 2620.0303 -> 2619.4243
$ws.Cells.Item(31, 13).Value = -2324.4243  # M31: -2325.0303 -> -2324.4243

# -- Row 34 --
$ws.Cells.Item(34, 8).Value = 3571.7317  # H34: 3572.2195 -> 3571.7317
$ws.Cells.Item(34, 9).Value = 2619.4243  # I34: 2620.0303 -> 2619.4243
$ws.Cells.Item(34, 11).Value = 2619.4243  # K34: 2620.0303 -> 2619.4243
$ws.Cells.Item(34, 13).Value = -2417.4243  # M34: -2418.0303 -> -2417.4243

# -- Row 58 --
$ws.Cells.Item(58, 8).Value = 2802.182  # H58: 2676.25 -> 2802.182
$ws.Cells.Item(58, 9).Value = 1457.5  # I58: 1424.2 -> 1457.5
$ws.Cells.Item(58, 11).Value = 1457.5  # K58: 1424.2 -> 1457.5
$ws.Cells.Item(58, 13).Value = -1254.5  # M58: -1221.2 -> -1254.5

# -- Row 62 --
$ws.Cells.Item(62, 8).Value = 20016120  # H62: 25017750 -> 20016120
$ws.Cells.Item(62, 9).Value = 25008150  # I62: 33341000 -> 25008150
$ws.Cells.Item(62, 11).Value = 25008150  # K62: 33341000 -> 25008150
$ws.Cells.Item(62, 13).Value = -25007526  # M62: -33340376 -> -25007526

# -- Row 65 --
$ws.Cells.Item(65, 8).Value = 20016120  # H65: 25017750 -> 20016120
$ws.Cells.Item(65, 9).Value = 25008150  # I65: 33341000 -> 25008150
$ws.Cells.Item(65, 11).Value = 125040750  # K65: 166705000 -> 125040750
$ws.Cells.Item(65, 13).Value = -125037630  # M65: -166701880 -> -125037630

# -- Row 134 --
$ws.Cells.Item(134, 8).Value = 2710.9473  # H134: 2789.4443 -> 2710.9473
$ws.Cells.Item(134, 9).Value = 2469.25  # I134: 2547.3333 -> 2469.25
$ws.Cells.Item(134, 11).Value = 7407.75  # K134: 7641.999899999999 -> 7407.75
$ws.Cells.Item(134, 13).Value = -4872.75  # M134: -5106.999899999999 -> -4872.75

# -- Row 136 --
$ws.Cells.Item(136, 8).Value = 2802.182  # H136: 2676.25 -> 2802.182
$ws.Cells.Item(136, 9).Value = 1457.5  # I136: 1424.2 -> 1457.5
$ws.Cells.Item(136, 11).Value = 4372.5  # K136: 4272.6 -> 4372.5
$ws.Cells.Item(136, 13).Value = -1822.5  # M136: -1722.6 -> -1822.5


# ======== Sheet: CUL ========
$ws = $wb.Worksheets.Item("CUL")

# -- Row 8 --
$ws.Cells.Item(8, 8).Value = 864.3333  # H8: 885.6667 -> 864.3333
$ws.Cells.Item(8, 9).Value = 864.3333  # I8: 885.6667 -> 864.3333
$ws.Cells.Item(8, 11).Value = 2592.9999  # K8: 2657.0001 -> 2592.9999
$ws.Cells.Item(8, 13).Value = -2453.9999  # M8: -2518.0001 -> -2453.9999

# -- Row 97 --
$ws.Cells.Item(97, 8).Value = 867916.7  # H97: 743926.5600000001 -> 867916.7
$ws.Cells.Item(97, 9).Value = 2500250  # I97: 1250371.5 -> 2500250
$ws.Cells.Item(97, 10).Value = 51750  # J97: 68666.664 -> 51750
$ws.Cells.Item(97, 11).Value = 7500750  # K97: 3751114.5 -> 7500750
$ws.Cells.Item(97, 12).Value = 155250  # L97: 205999.992 -> 155250
$ws.Cells.Item(97, 13).Value = -7500254  # M97: -3750618.5 -> -7500254
$ws.Cells.Item(97, 14).Value = -156242  # N97: -206991.992 -> -156242

# -- Row 113 --
$ws.Cells.Item(113, 8).Value = 1741.4  # H113: 1919.6428 -> 1741.4
$ws.Cells.Item(113, 9).Value = 825  # I113: 1000 -> 825
$ws.Cells.Item(113, 10).Value = 1882.3846  # J113: 1990.3846 -> 1882.3846
$ws.Cells.Item(113, 11).Value = 2475  # K113: 3000 -> 2475
$ws.Cells.Item(113, 12).Value = 5647.1538  # L113: 5971.1538 -> 5647.1538
$ws.Cells.Item(113, 13).Value = -305  # M113: -830 -> -305
$ws.Cells.Item(113, 14).Value = -9987.1538  # N113: -10311.1538 -> -9987.1538

# -- Row 114 --
$ws.Cells.Item(114, 8).Value = 896.375  # H114: 902.7778 -> 896.375
$ws.Cells.Item(114, 10).Value = 1513  # J114: 1401.2 -> 1513
$ws.Cells.Item(114, 12).Value = 4539  # L114: 4203.6 -> 4539
$ws.Cells.Item(114, 14).Value = -11047  # N114: -10711.6 -> -11047


# ======== Sheet: GSM ========
$ws = $wb.Worksheets.Item("GSM")

# -- Row 70 --
$ws.Cells.Item(70, 8).Value = 160384  # H70: 173207.75 -> 160384
$ws.Cells.Item(70, 9).Value = 403998.6  # I70: 503373.5 -> 403998.6
$ws.Cells.Item(70, 11).Value = 403998.6  # K70: 503373.5 -> 403998.6
$ws.Cells.Item(70, 13).Value = -403728.6  # M70: -503103.5 -> -403728.6

# -- Row 73 --
$ws.Cells.Item(73, 8).Value = 160384  # H73: 173207.75 -> 160384
$ws.Cells.Item(73, 9).Value = 403998.6  # I73: 503373.5 -> 403998.6
$ws.Cells.Item(73, 11).Value = 403998.6  # K73: 503373.5 -> 403998.6
$ws.Cells.Item(73, 13).Value = -403062.6  # M73: -502437.5 -> -403062.6

# -- Row 126 --
$ws.Cells.Item(126, 8).Value = 4777.8887  # H126: 5142.375 -> 4777.8887
$ws.Cells.Item(126, 9).Value = 2375.25  # I126: 2448.5715 -> 2375.25
$ws.Cells.Item(126, 11).Value = 7125.75  # K126: 7345.7145 -> 7125.75
$ws.Cells.Item(126, 13).Value = -4655.75  # M126: -4875.7145 -> -4655.75


# ======== Sheet: LTW ========
$ws = $wb.Worksheets.Item("LTW")

# -- Row 7 --
$ws.Cells.Item(7, 8).Value = 6631.25  # H7: 7133.2666 -> 6631.25
$ws.Cells.Item(7, 9).Value = 5508.5  # I7: 5833.4165 -> 5508.5
$ws.Cells.Item(7, 10).Value = 9999.5  # J7: 12332.667 -> 9999.5
$ws.Cells.Item(7, 11).Value = 5508.5  # K7: 5833.4165 -> 5508.5
$ws.Cells.Item(7, 12).Value = 9999.5  # L7: 12332.667 -> 9999.5
$ws.Cells.Item(7, 13).Value = -5396.5  # M7: -5721.4165 -> -5396.5
$ws.Cells.Item(7, 14).Value = -10223.5  # N7: -12556.667 -> -10223.5

# -- Row 62 --
$ws.Cells.Item(62, 8).Value = 45159  # H62: 0 -> 45159
$ws.Cells.Item(62, 10).Value = 45159  # J62: 0 -> 45159
$ws.Cells.Item(62, 12).Value = 45159  # L62: 0 -> 45159
$ws.Cells.Item(62, 14).Value = -46407  # N62: None -> -46407

# -- Row 65 --
$ws.Cells.Item(65, 8).Value = 45159  # H65: 0 -> 45159
$ws.Cells.Item(65, 10).Value = 45159  # J65: 0 -> 45159
$ws.Cells.Item(65, 12).Value = 135477  # L65: 0 -> 135477
$ws.Cells.Item(65, 14).Value = -141717  # N65: None -> -141717

# -- Row 100 --
$ws.Cells.Item(100, 8).Value = 4049.625  # H100: 3199.5715 -> 4049.625
$ws.Cells.Item(100, 9).Value = 2749.6667  # I100: 1299.6 -> 2749.6667
$ws.Cells.Item(100, 11).Value = 2749.6667  # K100: 1299.6 -> 2749.6667
$ws.Cells.Item(100, 13).Value = -2208.6667  # M100: -758.5999999999999 -> -2208.6667

# -- Row 122 --
$ws.Cells.Item(122, 8).Value = 4253.6665  # H122: 4703.9565 -> 4253.6665
$ws.Cells.Item(122, 9).Value = 4032.6538  # I122: 4463.227 -> 4032.6538
$ws.Cells.Item(122, 11).Value = 12097.9614  # K122: 13389.681 -> 12097.9614
$ws.Cells.Item(122, 13).Value = -9647.9614  # M122: -10939.681 -> -9647.9614

# -- Row 126 --
$ws.Cells.Item(126, 8).Value = 6631.25  # H126: 7133.2666 -> 6631.25
$ws.Cells.Item(126, 9).Value = 5508.5  # I126: 5833.4165 -> 5508.5
$ws.Cells.Item(126, 10).Value = 9999.5  # J126: 12332.667 -> 9999.5
$ws.Cells.Item(126, 11).Value = 16525.5  # K126: 17500.2495 -> 16525.5
$ws.Cells.Item(126, 12).Value = 29998.5  # L126: 36998.001 -> 29998.5
$ws.Cells.Item(126, 13).Value = -14055.5  # M126: -15030.2495 -> -14055.5
$ws.Cells.Item(126, 14).Value = -34938.5  # N126: -41938.001 -> -34938.5

# -- Row 132 --
$ws.Cells.Item(132, 8).Value = 3169.0217  # H132: 2937.4524 -> 3169.0217
$ws.Cells.Item(132, 9).Value = 2105.027  # I132: 2207.8057 -> 2105.027
$ws.Cells.Item(132, 10).Value = 7543.222  # J132: 7315.3335 -> 7543.222
$ws.Cells.Item(132, 11).Value = 6315.081  # K132: 6623.4171 -> 6315.081
$ws.Cells.Item(132, 12).Value = 22629.666  # L132: 21946.0005 -> 22629.666
$ws.Cells.Item(132, 13).Value = -3785.081  # M132: -4093.4171 -> -3785.081
$ws.Cells.Item(132, 14).Value = -27689.666  # N132: -27006.0005 -> -27689.666


# ======== Sheet: WVR ========
$ws = $wb.Worksheets.Item("WVR")

# -- Row 122 --
$ws.Cells.Item(122, 8).Value = 35716012  # H122: 35715950 -> 35716012
$ws.Cells.Item(122, 9).Value = 1956.6  # I122: 1945.3334 -> 1956.6
$ws.Cells.Item(122, 10).Value = 125001144  # J122: 250000000 -> 125001144
$ws.Cells.Item(122, 11).Value = 5869.799999999999  # K122: 5836.0002 -> 5869.799999999999
$ws.Cells.Item(122, 12).Value = 375003432  # L122: 750000000 -> 375003432
$ws.Cells.Item(122, 13).Value = -3419.799999999999  # M122: -3386.0002 -> -3419.799999999999
$ws.Cells.Item(122, 14).Value = -375008332  # N122: -750004900 -> -375008332

# -- Row 132 --
$ws.Cells.Item(132, 8).Value = 3715.7896  # H132: 4388.9414 -> 3715.7896
$ws.Cells.Item(132, 9).Value = 3450  # I132: 3840.8 -> 3450
$ws.Cells.Item(132, 11).Value = 10350  # K132: 11522.4 -> 10350
$ws.Cells.Item(132, 13).Value = -7820  # M132: -8992.400000000001 -> -7820

# -- Row 136 --
$ws.Cells.Item(136, 8).Value = 1787.4667  # H136: 1810.2307 -> 1787.4667
$ws.Cells.Item(136, 9).Value = 1794.0714  # I136: 1819.8334 -> 1794.0714
$ws.Cells.Item(136, 11).Value = 5382.2142  # K136: 5459.5002 -> 5382.2142
$ws.Cells.Item(136, 13).Value = -2832.2142  # M136: -2909.5002 -> -2832.2142

